# workflowStepSet base class for sets classes
# Adds a new "classes" worksheet (tracking which "set" related objects need
# adducts/setObjects/ionizedXXX/setThreshold/origFGNames/groupAlgorithm/
# analysisInfo) after the existing "components" sheet.

$wb = $excel.ActiveWorkbook

# --- add the new sheet at the end of the tab strip ---------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "classes"

# --- row labels / headers -----------------------------------------------
# (written in the original authoring order so new shared-string entries
# land at the same indices as the authored workbook)
$ws.Range("A2").Value = "featuresSet"
$ws.Range("B1").Value = "adducts"
$ws.Range("C1").Value = "setObjects"
$ws.Range("D1").Value = "ionizedXXX"
$ws.Range("A3").Value = "featuresGroupsSet"
$ws.Range("A4").Value = "MSPeakListsSet"
$ws.Range("A5").Value = "formulasSet"
$ws.Range("A6").Value = "compoundsSet"
$ws.Range("A7").Value = "componentsSet"
$ws.Range("E1").Value = "setThreshold"
$ws.Range("F1").Value = "origFGNames"
$ws.Range("G1").Value = "groupAlgorithm"
$ws.Range("H1").Value = "analysisInfo"

# --- centre-aligned marker grid -----------------------------------------
$ws.Range("B2:I7").HorizontalAlignment = -4108
$ws.Range("B8:D10").HorizontalAlignment = -4108

# --- "X" marks showing which columns apply to which class ---------------
$ws.Range("B2").Value = "X"
$ws.Range("C2").Value = "X"
$ws.Range("D2").Value = "X"

$ws.Range("G3").Value = "X"

$ws.Range("B4").Value = "X"
$ws.Range("C4").Value = "X"
$ws.Range("H4").Value = "X"

$ws.Range("B5").Value = "X"
$ws.Range("C5").Value = "X"
$ws.Range("E5").Value = "X"
$ws.Range("F5").Value = "X"

$ws.Range("B6").Value = "X"
$ws.Range("C6").Value = "X"
$ws.Range("E6").Value = "X"
$ws.Range("F6").Value = "X"

$ws.Range("C7").Value = "X"

# --- column widths / fit --------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 17.85546875
$ws.Columns.Item(4).ColumnWidth = 11
$ws.Columns.Item(5).ColumnWidth = 12.5703125
$ws.Columns.Item(6).ColumnWidth = 12.85546875
$ws.Columns.Item(7).ColumnWidth = 15.140625
$ws.Columns.Item(8).ColumnWidth = 11.5703125

# --- selection on the new sheet matches the authored file ----------------
$ws.Range("C3").Select()

# make sure the new sheet is the active / visible tab (matches tabSelected
# moving from "components" to "classes")
$ws.Activate()
